$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update E12's text to append a trailing newline after "Đã có giao diện "
$ws.Range("E12").Value = "Đã có giao diện " + [char]10

# Fill in E13 with the new status text and F13 with 50% completion
$ws.Range("E13").Value = "Xây dựng cơ bản trên Cisco" + [char]10 + "Bổ sung ERD" + [char]10 + "Bổ sung Flowchart" + [char]10 + "Chỉnh sửa thiết kế chức năng" + [char]10 + "Code database theo sơ đồ ERD"

# Apply the percentage number format (matching the other weeks' % column) and set the value
$ws.Range("F13").NumberFormat = "0%"
$ws.Range("F13").Value = 0.5

# Adjust row 13 height to fit the newly added multi-line text
$ws.Rows.Item(13).RowHeight = 70.8

# Update the sheet view's scroll position and active selection
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Range("H13").Select()
